$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" for the
# 416f48b4 file row (row 2) to reflect the newly generated handback report.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-03 04:53:11"

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the 416f48b4 file row (row 2).
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-03 04:53:05"
$wsZhCn.Range("K2").Value = "2016-09-03 04:53:31"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback
# DateTime for the 416f48b4 file row (row 2).
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-03 04:53:11"
$wsDeDe.Range("K2").Value = "2016-09-03 04:53:38"
